# Natmi following Dr Hou advice
#
# The LR-pair grid (natmiOut/YoungD0/LR-pairs_lrc2p/Cadm1-Cadm1.xlsx) is re-run with an
# additional "FAPs" cell cluster alongside the pre-existing "ECs" and "sCs" clusters, so
# the Sheet1 table grows from a 2x2 (4-row) sending/target cluster grid to a full 3x3
# (9-row) cross-product grid, with refreshed NATMI specificity/expression statistics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster=ECs, Target cluster=ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cadm1"
$ws.Range("C2").Value = "Cadm1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.716522666666667
$ws.Range("H2").Value = 14.149568
$ws.Range("I2").Value = 0.530211572117814
$ws.Range("J2").Value = 0.530211572117814
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.716522666666667
$ws.Range("N2").Value = 14.149568
$ws.Range("O2").Value = 0.530211572117814
$ws.Range("P2").Value = 0.530211572117814
$ws.Range("Q2").Value = 22.24558606518045
$ws.Range("R2").Value = 200.210274586624
$ws.Range("S2").Value = 0.2811243112076439
$ws.Range("T2").Value = 0.2811243112076439

# Row 3: Sending cluster=ECs, Target cluster=FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cadm1"
$ws.Range("C3").Value = "Cadm1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.716522666666667
$ws.Range("H3").Value = 14.149568
$ws.Range("I3").Value = 0.530211572117814
$ws.Range("J3").Value = 0.530211572117814
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.236474
$ws.Range("N3").Value = 0.709422
$ws.Range("O3").Value = 0.02658340904223816
$ws.Range("P3").Value = 0.02658340904223817
$ws.Range("Q3").Value = 1.115334981077333
$ws.Range("R3").Value = 10.038014829696
$ws.Range("S3").Value = 0.01409483110053601
$ws.Range("T3").Value = 0.01409483110053601

# Row 4: Sending cluster=ECs, Target cluster=sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cadm1"
$ws.Range("C4").Value = "Cadm1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.716522666666667
$ws.Range("H4").Value = 14.149568
$ws.Range("I4").Value = 0.530211572117814
$ws.Range("J4").Value = 0.530211572117814
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.942551666666667
$ws.Range("N4").Value = 11.827655
$ws.Range("O4").Value = 0.4432050188399478
$ws.Range("P4").Value = 0.4432050188399479
$ws.Range("Q4").Value = 18.59513430033778
$ws.Range("R4").Value = 167.35620870304
$ws.Range("S4").Value = 0.2349924298096341
$ws.Range("T4").Value = 0.2349924298096341

# Row 5: Sending cluster=FAPs, Target cluster=ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cadm1"
$ws.Range("C5").Value = "Cadm1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.236474
$ws.Range("H5").Value = 0.709422
$ws.Range("I5").Value = 0.02658340904223816
$ws.Range("J5").Value = 0.02658340904223817
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.716522666666667
$ws.Range("N5").Value = 14.149568
$ws.Range("O5").Value = 0.530211572117814
$ws.Range("P5").Value = 0.530211572117814
$ws.Range("Q5").Value = 1.115334981077333
$ws.Range("R5").Value = 10.038014829696
$ws.Range("S5").Value = 0.01409483110053601
$ws.Range("T5").Value = 0.01409483110053601

# Row 6: Sending cluster=FAPs, Target cluster=FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cadm1"
$ws.Range("C6").Value = "Cadm1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.236474
$ws.Range("H6").Value = 0.709422
$ws.Range("I6").Value = 0.02658340904223816
$ws.Range("J6").Value = 0.02658340904223817
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.236474
$ws.Range("N6").Value = 0.709422
$ws.Range("O6").Value = 0.02658340904223816
$ws.Range("P6").Value = 0.02658340904223817
$ws.Range("Q6").Value = 0.05591995267599999
$ws.Range("R6").Value = 0.503279574084
$ws.Range("S6").Value = 0.0007066776363069497
$ws.Range("T6").Value = 0.0007066776363069501

# Row 7: Sending cluster=FAPs, Target cluster=sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cadm1"
$ws.Range("C7").Value = "Cadm1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.236474
$ws.Range("H7").Value = 0.709422
$ws.Range("I7").Value = 0.02658340904223816
$ws.Range("J7").Value = 0.02658340904223817
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.942551666666667
$ws.Range("N7").Value = 11.827655
$ws.Range("O7").Value = 0.4432050188399478
$ws.Range("P7").Value = 0.4432050188399479
$ws.Range("Q7").Value = 0.9323109628233333
$ws.Range("R7").Value = 8.390798665409999
$ws.Range("S7").Value = 0.0117819003053952
$ws.Range("T7").Value = 0.01178190030539521

# Row 8: Sending cluster=sCs, Target cluster=ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cadm1"
$ws.Range("C8").Value = "Cadm1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.942551666666667
$ws.Range("H8").Value = 11.827655
$ws.Range("I8").Value = 0.4432050188399478
$ws.Range("J8").Value = 0.4432050188399479
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.716522666666667
$ws.Range("N8").Value = 14.149568
$ws.Range("O8").Value = 0.530211572117814
$ws.Range("P8").Value = 0.530211572117814
$ws.Range("Q8").Value = 18.59513430033778
$ws.Range("R8").Value = 167.35620870304
$ws.Range("S8").Value = 0.2349924298096341
$ws.Range("T8").Value = 0.2349924298096341

# Row 9: Sending cluster=sCs, Target cluster=FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cadm1"
$ws.Range("C9").Value = "Cadm1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.942551666666667
$ws.Range("H9").Value = 11.827655
$ws.Range("I9").Value = 0.4432050188399478
$ws.Range("J9").Value = 0.4432050188399479
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.236474
$ws.Range("N9").Value = 0.709422
$ws.Range("O9").Value = 0.02658340904223816
$ws.Range("P9").Value = 0.02658340904223817
$ws.Range("Q9").Value = 0.9323109628233333
$ws.Range("R9").Value = 8.390798665409999
$ws.Range("S9").Value = 0.0117819003053952
$ws.Range("T9").Value = 0.01178190030539521

# Row 10: Sending cluster=sCs, Target cluster=sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cadm1"
$ws.Range("C10").Value = "Cadm1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.942551666666667
$ws.Range("H10").Value = 11.827655
$ws.Range("I10").Value = 0.4432050188399478
$ws.Range("J10").Value = 0.4432050188399479
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.942551666666667
$ws.Range("N10").Value = 11.827655
$ws.Range("O10").Value = 0.4432050188399478
$ws.Range("P10").Value = 0.4432050188399479
$ws.Range("Q10").Value = 15.54371364433611
$ws.Range("R10").Value = 139.893422799025
$ws.Range("S10").Value = 0.1964306887249185
$ws.Range("T10").Value = 0.1964306887249186

